# Smart fan BOM update: found the fan header part at Digikey (was previously
# an unresolved Ebay listing with a "need to find..." comment).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Point the existing hyperlink on B6 to the Digikey product page instead of
# the old Ebay listing, and drop the old Ebay-specific anchor/location.
foreach ($h in @($ws.Hyperlinks)) {
    if ($h.Address -like "*ebay*") {
        $h.Address = "https://www.digikey.com/product-detail/en/molex-llc/0022232061/WM4212-ND/26691"
        $h.SubAddress = ""
    }
}

# Update row 6 with the newly found part: name, vendor, price and quantity.
$ws.Range("A6").Value = "fan header"
$ws.Range("B6").Value = "Digikey"
$ws.Range("C6").Value = 0.46
$ws.Range("D6").Value = 3

# The "need to find..." comment is no longer needed now that the part's found.
$ws.Range("F6").ClearContents()
